$wb = $excel.ActiveWorkbook
$docentes = $wb.Worksheets.Item("docentes")
$ws = $wb.Worksheets.Item("asesorias")

# Shift old header values (B1:F1) one column to the right (C1:G1) before overwriting B1/C1
$ws.Range("G1").Value = $ws.Range("F1").Value2
$ws.Range("F1").Value = $ws.Range("E1").Value2
$ws.Range("E1").Value = $ws.Range("D1").Value2
$ws.Range("D1").Value = $ws.Range("C1").Value2
$ws.Range("C1").Value = $ws.Range("B1").Value2

# New header text (C1 first so shared-string order matches: ...,"Usuario Asesor","Usuario Estudiante")
$ws.Range("C1").Value = "Usuario Asesor"
$ws.Range("B1").Value = "Usuario Estudiante"

# Give C1 the bold/themed (no border) style, matching the "Oficina" header style used in docentes!L1
$docentes.Range("L1").Copy()
$ws.Range("C1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Remove the old sample data row
$ws.Rows.Item(2).Delete()

# Make asesorias the active sheet/tab with the new selection
$ws.Activate() | Out-Null
$ws.Range("A2:B2").Select() | Out-Null
Write-Host "done"
